$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the Send_email / AI_Summary config value: update the email address used by Config
$ws.Range("B7").Value = "tunaru.alexandra2005@gmail.com"

# Excel re-autofits the row height for row 7 after the longer text is entered
$ws.Rows.Item(7).RowHeight = 15

# Reflect the last active selection recorded in the saved file
$ws.Range("D6").Select()
